function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "41.854.03"
    "E2" = "  +5.42%  "
    "D3" = "2.228.31"
    "E3" = "  +2.38%  "
    "D4" = "0.999"
    "E4" = "  -0.18%  "
    "D5" = "231.59"
    "E5" = "  +2.04%  "
    "E6" = "  -0.65%  "
    "D7" = "61.55"
    "E7" = "  -2.58%  "
    "E8" = "  -0.04%  "
    "D9" = "0.404"
    "E9" = "  +2.94%  "
    "D10" = "59.10"
    "E10" = "  +1.02%  "
    "D11" = "0.0897"
    "E11" = "  +5.23%  "
    "E12" = "  -0.14%  "
    "D13" = "2.561.36"
    "E13" = "  +2.51%  "
    "D14" = "15.64"
    "E14" = "  -1.88%  "
    "D15" = "22.03"
    "E15" = "  +0.76%  "
    "D16" = "0.800"
    "E16" = "  -1.64%  "
    "E17" = "  +1.21%  "
    "D18" = "2.257.25"
    "E18" = "  +3.69%  "
    "D19" = "41.802.13"
    "E19" = "  +5.27%  "
    "B20" = "Litecoin"
    "C20" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D20" = "72.07"
    "E20" = "  +0.35%  "
    "B21" = "ShibaInu"
    "C21" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D21" = "0.0₃0900"
    "E21" = "  -2.26%  "
    "D22" = "6.03"
    "E22" = "  +0.48%  "
    "D23" = "249.31"
    "E23" = "  +8.49%  "
    "E25" = "  +1.87%  "
    "E26" = "  -0.63%  "
    "D27" = "9.60"
    "E27" = "  -0.18%  "
    "B28" = "Monero"
    "C28" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D28" = "167.90"
    "E28" = "  -1.81%  "
    "B29" = "Kaspa"
    "C29" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D29" = "0.142"
    "E29" = "  +1.07%  "
    "D30" = "20.01"
    "E30" = "  +0.72%  "
    "E31" = "  -3.06%  "
    "D32" = "2.71"
    "E32" = "  +1.21%  "
    "E34" = "  +7.64%  "
    "D35" = "4.66"
    "E35" = "  +2.78%  "
    "D36" = "0.0637"
    "E36" = "  +3.08%  "
    "D37" = "6.62"
    "E37" = "  -5.01%  "
    "E38" = "  -6.42%  "
    "D39" = "2.36"
    "E39" = "  -1.78%  "
    "D40" = "0.000267"
    "E40" = "  +37.90%  "
    "D41" = "0.997"
    "E41" = "  -0.40%  "
    "B42" = "VeChain"
    "C42" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D42" = "0.0239"
    "E42" = "  +4.55%  "
    "B43" = "FTXToken"
    "C43" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D43" = "4.85"
    "E43" = "  -2.69%  "
    "D44" = "8.57"
    "E44" = "  +8.45%  "
    "D45" = "1.22"
    "E45" = "  -0.24%  "
    "D46" = "0.0972"
    "E46" = "  +5.56%  "
    "D47" = "98.93"
    "E47" = "  -3.56%  "
    "D48" = "1.477.77"
    "E48" = "  -2.32%  "
    "B49" = "HuobiToken"
    "C49" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D49" = "2.81"
    "E49" = "  +0.24%  "
    "B50" = "InjectiveProtocol"
    "C50" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D50" = "16.46"
    "E50" = "  -7.29%  "
    "D51" = "52.52"
    "E51" = "  +5.69%  "
}

foreach ($addr in $updates.Keys) {
    Set-TextValue $ws $addr $updates[$addr]
}
